$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: 16/09/2025, JM Empilhadeira job vacancy ---
$ws.Range("A16").Value = 45916
$ws.Range("A16").NumberFormat = "dd/mm/yyyy"
$ws.Range("C16").Value = "09:00-16:00"
$ws.Range("D16").Value = "Vagas - Auxiliar de Logística - Auxilia de Armazém "
$ws.Range("D16").WrapText = $true
$ws.Range("E16").Value = "JM Empilhadeira "
$ws.Rows(16).RowHeight = 14.25

# --- Row 17: 17/09/2025, Desenvolve Já job vacancy ---
$ws.Range("A17").Value = 45917
$ws.Range("A17").NumberFormat = "dd/mm/yyyy"
$ws.Range("C17").Value = "09:00-16:00"
$ws.Range("D17").Value = "Vagas - Telemarketing "
$ws.Range("D17").WrapText = $true
$ws.Range("E17").Value = "Desenvolve Já "
$ws.Rows(17).RowHeight = 14.25

# --- Row 18: 18/09/2025, Líderes Consultoria job vacancies (longer text -> taller row) ---
$ws.Range("A18").Value = 45918
$ws.Range("A18").NumberFormat = "dd/mm/yyyy"
$ws.Range("C18").Value = "09:00-16:00"
$ws.Range("D18").Value = "Vagas - 7 vagas de vendedor(a) ótica/2 recepcionista para cliente/2 Gerentes de vendas (ótica)"
$ws.Range("D18").WrapText = $true
$ws.Range("E18").Value = "Líderes Consultoria e Treinamentos"
$ws.Rows(18).RowHeight = 28.5

# --- Row 19: 19/09/2025, Grau Técnico job vacancies ---
$ws.Range("A19").Value = 45919
$ws.Range("A19").NumberFormat = "dd/mm/yyyy"
$ws.Range("C19").Value = "09:00-16:00"
$ws.Range("D19").Value = "Vagas - Vendedor Interno /Vendedor externo /Estágio ( Setor Financeiro)"
$ws.Range("E19").Value = "Grau Técnico"
$ws.Rows(19).RowHeight = 14.25

# --- Row 20: 18/09/2025, Manhã - Emissao CIN ---
$ws.Range("A20").Value = 45918
$ws.Range("A20").NumberFormat = "dd/mm/yyyy"
$ws.Range("B20").Value = "Manhã"
$ws.Range("C20").Value = "08:00-12:00"
$ws.Range("D20").Value = "Emissão de novas CIN (Carteira de Identidade Nacional)."
$ws.Range("E20").Value = "Instituto de Cidadania Digital Félix Pacheco"
$ws.Rows(20).RowHeight = 14.25

# --- Row 21: 19/09/2025, Manhã - Emissao CIN ---
$ws.Range("A21").Value = 45919
$ws.Range("A21").NumberFormat = "dd/mm/yyyy"
$ws.Range("B21").Value = "Manhã"
$ws.Range("C21").Value = "08:00-12:00"
$ws.Range("C21").NumberFormat = "h:mm"
$ws.Range("D21").Value = "Emissão de novas CIN (Carteira de Identidade Nacional)."
$ws.Range("D21").WrapText = $true
$ws.Range("D21").WrapText = $false
$ws.Range("E21").Value = "Instituto de Cidadania Digital Félix Pacheco"
$ws.Range("E21").WrapText = $true
$ws.Range("E21").WrapText = $false
$ws.Rows(21).RowHeight = 14.25

# --- Column widths (D widened/custom, E widened) ---
$ws.Columns("D").ColumnWidth = 63
$ws.Columns("E").ColumnWidth = 37
